# repull data, push all data, mean calculation
# Update the "dSF" column (F) with freshly re-pulled delta-to-final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 9
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 0
